# Rewrite the "KEY ACHIEVEMENTS AND IMPACT" > "Impact" bullet list so that
# it contains concise, impact-focused accomplishment statements instead of
# the previous job-duty style bullets.

$d = $word.ActiveDocument

# Locate the "KEY ACHIEVEMENTS AND IMPACT" heading paragraph so we can work
# relative to it (the bullet text under "Partner - Siege Analytics" higher
# up in the document is very similar, so we must not rely on plain text
# search across the whole document).
$achievementsHeadingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq "KEY ACHIEVEMENTS AND IMPACT") {
        $achievementsHeadingIndex = $i
        break
    }
}

if ($achievementsHeadingIndex -eq -1) {
    throw "Could not find KEY ACHIEVEMENTS AND IMPACT heading"
}

# The bullet list begins two paragraphs after the heading ("Impact" sub-
# heading sits in between).
$firstBulletIndex = $achievementsHeadingIndex + 2

# Update the first four bullets with the new, tightened copy. (No trailing
# `r -- Range.Text already replaces just the paragraph's own text and
# keeps its existing paragraph mark; appending `r would insert an extra
# paragraph break.)
$d.Paragraphs.Item($firstBulletIndex + 0).Range.Text = "• Predictive excellence: Achieved 87% voter turnout accuracy vs. 71% industry standard"
$d.Paragraphs.Item($firstBulletIndex + 1).Range.Text = "• Reduced polling margins from ±4.2% to ±2.1%"
$d.Paragraphs.Item($firstBulletIndex + 2).Range.Text = "• Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis"
$d.Paragraphs.Item($firstBulletIndex + 3).Range.Text = "• Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%"

# Remove the trailing two bullets entirely (they are dropped in the new
# version of the achievements list).
$d.Paragraphs.Item($firstBulletIndex + 5).Range.Delete()
$d.Paragraphs.Item($firstBulletIndex + 4).Range.Delete()
